$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.406.66"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.14"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9996"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "284.77"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3620"
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.43"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3336"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07392"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.76"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.932"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.889"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.564.52"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001105"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.13"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06688"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.11"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.408.35"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.415"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.558"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.40"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.000"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.28"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.737.89"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.136"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.994"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.796"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08281"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02403"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06379"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2210"
$ws.Range("E39").Value = "  -3.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.287"
$ws.Range("E40").Value = "  -6.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.318"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.16"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6079"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.755"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5768"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.51"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.213"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07204"
$ws.Range("E51").Value = "  -1.47%  "
